# Data retrieved - Sat Jul  3 18:24:40 UTC 2021
#
# Corrects the timestamp recorded for row 66 (tiny floating-point fix to the
# serial date/time value) and appends a newly-retrieved data row (row 67).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix up the previously-recorded timestamp in A66 (same calendar day, just a
# corrected fractional-second value).
$ws.Cells.Item(66, 1).Value = 44379.76524352431

# Append the newly retrieved row of job numbers.
$ws.Cells.Item(67, 1).Value = 44380.76713137147
$ws.Cells.Item(67, 2).Value = 78860
$ws.Cells.Item(67, 3).Value = 66459
$ws.Cells.Item(67, 4).Value = 3626
$ws.Cells.Item(67, 5).Value = 2133
$ws.Cells.Item(67, 6).Value = 1510
$ws.Cells.Item(67, 7).Value = 20932
$ws.Cells.Item(67, 8).Value = 1530
$ws.Cells.Item(67, 9).Value = 873
$ws.Cells.Item(67, 10).Value = 192
